$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 232.38889
$ws.Range("I53").Value = 234.33333
$ws.Range("J53").Value = 228.5
$ws.Range("K53").Value = 234.33333
$ws.Range("L53").Value = 228.5
$ws.Range("M53").Value = 402.66667
$ws.Range("N53").Value = -1502.5
$ws.Range("H64").Value = 3759.8
$ws.Range("I64").Value = 5999
$ws.Range("K64").Value = 5999
$ws.Range("M64").Value = -5751
$ws.Range("H67").Value = 3759.8
$ws.Range("I67").Value = 5999
$ws.Range("K67").Value = 5999
$ws.Range("M67").Value = -5141
$ws.Range("H86").Value = 3883.3333
$ws.Range("I86").Value = 1598
$ws.Range("J86").Value = 4714.364
$ws.Range("K86").Value = 1598
$ws.Range("L86").Value = 4714.364
$ws.Range("M86").Value = -475
$ws.Range("N86").Value = -6960.364
$ws.Range("H89").Value = 3883.3333
$ws.Range("I89").Value = 1598
$ws.Range("J89").Value = 4714.364
$ws.Range("K89").Value = 7990
$ws.Range("L89").Value = 23571.82
$ws.Range("M89").Value = -2374
$ws.Range("N89").Value = -34803.82
$ws.Range("H131").Value = 4214.875
$ws.Range("I131").Value = 1142.1428
$ws.Range("K131").Value = 3426.4284
$ws.Range("M131").Value = 1613.5716
$ws.Range("H132").Value = 1985.9546
$ws.Range("I132").Value = 1913.8572
$ws.Range("K132").Value = 5741.571599999999
$ws.Range("M132").Value = -3211.571599999999
$ws.Range("H137").Value = 2143.3333
$ws.Range("I137").Value = 1065.3572
$ws.Range("J137").Value = 4299.2856
$ws.Range("K137").Value = 3196.0716
$ws.Range("L137").Value = 12897.8568
$ws.Range("M137").Value = -646.0715999999998
$ws.Range("N137").Value = -17997.8568
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1384.2206
$ws.Range("I61").Value = 1384.2206
$ws.Range("K61").Value = 1384.2206
$ws.Range("M61").Value = -1172.2206
$ws.Range("H74").Value = 2984.8333
$ws.Range("I74").Value = 1014.5
$ws.Range("J74").Value = 5447.75
$ws.Range("K74").Value = 1014.5
$ws.Range("L74").Value = 5447.75
$ws.Range("M74").Value = -140.5
$ws.Range("N74").Value = -7195.75
$ws.Range("H77").Value = 2984.8333
$ws.Range("I77").Value = 1014.5
$ws.Range("J77").Value = 5447.75
$ws.Range("K77").Value = 5072.5
$ws.Range("L77").Value = 27238.75
$ws.Range("M77").Value = -704.5
$ws.Range("N77").Value = -35974.75
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H136").Value = 1384.2206
$ws.Range("I136").Value = 1384.2206
$ws.Range("K136").Value = 4152.6618
$ws.Range("M136").Value = -1602.6618

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2335.25
$ws.Range("I134").Value = 2115.56
$ws.Range("K134").Value = 6346.68
$ws.Range("M134").Value = -3811.68

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3546.3684
$ws.Range("I58").Value = 2453.4285
$ws.Range("J58").Value = 4183.9165
$ws.Range("K58").Value = 2453.4285
$ws.Range("L58").Value = 4183.9165
$ws.Range("M58").Value = -2250.4285
$ws.Range("N58").Value = -4589.9165
$ws.Range("H86").Value = 12232.375
$ws.Range("I86").Value = 9365.75
$ws.Range("K86").Value = 9365.75
$ws.Range("M86").Value = -8242.75
$ws.Range("H89").Value = 12232.375
$ws.Range("I89").Value = 9365.75
$ws.Range("K89").Value = 46828.75
$ws.Range("M89").Value = -41212.75
$ws.Range("H105").Value = 3324.8
$ws.Range("I105").Value = 958.3333
$ws.Range("J105").Value = 6874.5
$ws.Range("K105").Value = 958.3333
$ws.Range("L105").Value = 6874.5
$ws.Range("M105").Value = 788.6667
$ws.Range("N105").Value = -10368.5
$ws.Range("H107").Value = 180
$ws.Range("I107").Value = 180
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 180
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1740
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 3709.158
$ws.Range("I134").Value = 3177.923
$ws.Range("J134").Value = 4860.1665
$ws.Range("K134").Value = 9533.769
$ws.Range("L134").Value = 14580.4995
$ws.Range("M134").Value = -6998.769
$ws.Range("N134").Value = -19650.4995
$ws.Range("H136").Value = 3546.3684
$ws.Range("I136").Value = 2453.4285
$ws.Range("J136").Value = 4183.9165
$ws.Range("K136").Value = 7360.2855
$ws.Range("L136").Value = 12551.7495
$ws.Range("M136").Value = -4810.2855
$ws.Range("N136").Value = -17651.7495
$ws.Range("H141").Value = 168331.33
$ws.Range("J141").Value = 168331.33
$ws.Range("L141").Value = 168331.33
$ws.Range("N141").Value = -178691.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3300
$ws.Range("J5").Value = 1600
$ws.Range("L5").Value = 4800
$ws.Range("N5").Value = -5024
$ws.Range("H44").Value = 1539.92
$ws.Range("I44").Value = 2000
$ws.Range("J44").Value = 1499.9131
$ws.Range("K44").Value = 6000
$ws.Range("L44").Value = 4499.7393
$ws.Range("M44").Value = -5602
$ws.Range("N44").Value = -5295.7393
$ws.Range("H121").Value = 1593.5
$ws.Range("I121").Value = 369.6
$ws.Range("K121").Value = 1108.8
$ws.Range("M121").Value = 201.1999999999998
$ws.Range("H122").Value = 1595.3334
$ws.Range("I122").Value = 1893
$ws.Range("K122").Value = 17037
$ws.Range("M122").Value = -14587
$ws.Range("H131").Value = 4866.1113
$ws.Range("I131").Value = 6081.6665
$ws.Range("J131").Value = 4258.3335
$ws.Range("K131").Value = 18244.9995
$ws.Range("L131").Value = 12775.0005
$ws.Range("M131").Value = -13204.9995
$ws.Range("N131").Value = -22855.0005
$ws.Range("H132").Value = 2886.111
$ws.Range("J132").Value = 3757.6
$ws.Range("L132").Value = 33818.4
$ws.Range("N132").Value = -38878.4
$ws.Range("H135").Value = 3300
$ws.Range("J135").Value = 1600
$ws.Range("L135").Value = 14400
$ws.Range("N135").Value = -19470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4166.3335
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 4166.3335
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H102").Value = 497.5
$ws.Range("I102").Value = 497.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 497.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 1124.5
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 41100.11
$ws.Range("I122").Value = 3989.7856
$ws.Range("J122").Value = 81065.08
$ws.Range("K122").Value = 11969.3568
$ws.Range("L122").Value = 243195.24
$ws.Range("M122").Value = -9519.356800000001
$ws.Range("N122").Value = -248095.24
$ws.Range("H132").Value = 5857
$ws.Range("I132").Value = 8400
$ws.Range("J132").Value = 5221.25
$ws.Range("K132").Value = 25200
$ws.Range("L132").Value = 15663.75
$ws.Range("M132").Value = -22670
$ws.Range("N132").Value = -20723.75
$ws.Range("H141").Value = 58259.6
$ws.Range("J141").Value = 58259.6
$ws.Range("L141").Value = 58259.6
$ws.Range("N141").Value = -68619.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 999.0909
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H93").Value = 613.5
$ws.Range("I93").Value = 599.5
$ws.Range("J93").Value = 627.5
$ws.Range("K93").Value = 599.5
$ws.Range("L93").Value = 627.5
$ws.Range("M93").Value = 648.5
$ws.Range("N93").Value = -3123.5
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 5044.1113
$ws.Range("I122").Value = 3678.8
$ws.Range("J122").Value = 5569.231
$ws.Range("K122").Value = 11036.4
$ws.Range("L122").Value = 16707.693
$ws.Range("M122").Value = -8586.400000000001
$ws.Range("N122").Value = -21607.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10400
$ws.Range("H104").Value = 22846.5
$ws.Range("J104").Value = 22846.5
$ws.Range("L104").Value = 22846.5
$ws.Range("N104").Value = -29834.5
$ws.Range("H107").Value = 806.1579
$ws.Range("I107").Value = 430
$ws.Range("J107").Value = 1859.4
$ws.Range("K107").Value = 1290
$ws.Range("L107").Value = 5578.200000000001
$ws.Range("M107").Value = 630
$ws.Range("N107").Value = -9418.200000000001
$ws.Range("H122").Value = 620.6111
$ws.Range("I122").Value = 573.5
$ws.Range("J122").Value = 997.5
$ws.Range("K122").Value = 1720.5
$ws.Range("L122").Value = 2992.5
$ws.Range("M122").Value = 729.5
$ws.Range("N122").Value = -7892.5
$ws.Range("H126").Value = 2468.4443
$ws.Range("I126").Value = 871.6667
$ws.Range("K126").Value = 2615.0001
$ws.Range("M126").Value = -145.0001000000002
